# Double every raw count in the "totals" and "cleavages" sheets (data range
# B2:U21). The "probs" sheet is left untouched; it consists purely of
# formulas (cleavages/totals ratios) and recalculates automatically.

$wb = $excel.ActiveWorkbook

$sheetNames = @("totals", "cleavages")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 21; $row++) {
        for ($col = 2; $col -le 21; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $val = $cell.Value2
            if ($val -ne $null) {
                $cell.Value = ($val * 2)
            }
        }
    }
}
